$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "44"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "18"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "5"

# Row 3
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "42"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2"

# Row 4
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "24"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "23"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "1"

# Row 5
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "7"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0"

# Row 7
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "8"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "12"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "0"
